# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as a new data row at row 676,
# pushing the existing rows 676:772 down to 677:773 (dimension grows to A1:R773).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 676 - this shifts rows 676:772 down
# to 677:773, carrying their values/formatting with them (matches the diff,
# where every existing row below 675 simply moves down by one row).
$ws.Rows.Item(676).Insert()

# Populate the newly inserted row 676 with the new observation.
$ws.Range("A676").Value2 = 3
$ws.Range("B676").Value2 = "Femacal de La Calera"
$ws.Range("C676").Value2 = "Coquimbo"
$ws.Range("D676").Value2 = 45077
$ws.Range("E676").Value2 = 5
$ws.Range("F676").Value2 = 100112032
$ws.Range("G676").Value2 = "Zapallo italiano"
$ws.Range("H676").Value2 = "Sin especificar"
$ws.Range("I676").Value2 = "Primera"
$ws.Range("J676").Value2 = 262
$ws.Range("K676").Value2 = 8500
$ws.Range("L676").Value2 = 9000
$ws.Range("M676").Value2 = 8783
$ws.Range("N676").Value2 = "`$/caja 60 unidades"
$ws.Range("O676").Value2 = "Región de Arica y Parinacota"
$ws.Range("P676").Value2 = 146
$ws.Range("Q676").Value2 = 60
$ws.Range("R676").Value2 = "Hortaliza"

# Keep the same date number format used by the rest of column D.
$ws.Range("D676").NumberFormat = $ws.Range("D677").NumberFormat
